$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.367.20'
$ws.Range("E2").Value = '  -3.62%  '

$ws.Range("D3").Value = '3.124.51'
$ws.Range("E3").Value = '  -4.63%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.75'
$ws.Range("E5").Value = '  -4.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.45'
$ws.Range("E6").Value = '  -9.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("E8").Value = '  -9.07%  '

$ws.Range("D9").Value = '3.115.43'
$ws.Range("E9").Value = '  -4.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  -2.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.115'
$ws.Range("E11").Value = '  -7.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.377'
$ws.Range("E12").Value = '  -5.23%  '

$ws.Range("D13").Value = '3.661.51'
$ws.Range("E13").Value = '  -5.15%  '

$ws.Range("E14").Value = '  -1.27%  '

$ws.Range("D15").Value = '63.368.92'
$ws.Range("E15").Value = '  -3.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.69'
$ws.Range("E16").Value = '  -6.35%  '

$ws.Range("D17").Value = '3.119.24'
$ws.Range("E17").Value = '  -4.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000152'
$ws.Range("E18").Value = '  -6.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '397.86'
$ws.Range("E19").Value = '  -5.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.19'
$ws.Range("E20").Value = '  -4.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.41'
$ws.Range("E21").Value = '  -4.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.02'
$ws.Range("E22").Value = '  -3.04%  '

$ws.Range("E23").Value = '  +0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.14'
$ws.Range("E24").Value = '  -5.38%  '

$ws.Range("E25").Value = '  -3.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.475'
$ws.Range("E26").Value = '  -6.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000100'
$ws.Range("E27").Value = '  -11.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.67'
$ws.Range("E28").Value = '  -7.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.78'
$ws.Range("E31").Value = '  -7.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.81'
$ws.Range("E32").Value = '  -6.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.22'
$ws.Range("E33").Value = '  -5.05%  '

$ws.Range("E34").Value = '  -6.82%  '

$ws.Range("E35").Value = '  -7.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '151.87'
$ws.Range("E36").Value = '  -3.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.32'
$ws.Range("E37").Value = '  -8.02%  '

$ws.Range("D38").Value = '2.742.37'
$ws.Range("E38").Value = '  -3.74%  '

$ws.Range("E39").Value = '  -8.10%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.27'
$ws.Range("E40").Value = '  -11.06%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").Value = '  -6.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.21'
$ws.Range("E42").Value = '  -3.36%  '

$ws.Range("E43").Value = '  -7.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0610'
$ws.Range("E44").Value = '  -3.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.40'
$ws.Range("E45").Value = '  -8.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0256'
$ws.Range("E46").Value = '  -4.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.67'
$ws.Range("E47").Value = '  -9.15%  '

$ws.Range("E48").Value = '  -0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '278.98'
$ws.Range("E49").Value = '  -10.28%  '

$ws.Range("E50").Value = '  -4.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.48'
$ws.Range("E51").Value = '  +0.94%  '
